# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column headers in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting from the adjacent existing header cell (AC1)
# so the new header cells share the same bold/centered/bordered style.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Data rows: every row (2 through 55) gets the same team record values.
$lastRow = 55
$ws.Range("AD2:AD" + $lastRow).Value = 90
$ws.Range("AE2:AE" + $lastRow).Value = 73
$ws.Range("AF2:AF" + $lastRow).Value = 0

$excel.CutCopyMode = 0
